# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 ("R") gets updated with Wild Card round stats added in
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 396
$wsOff.Range("C3").Value = 268
$wsOff.Range("D3").Value = 70
$wsOff.Range("E3").Value = 24
$wsOff.Range("F3").Value = 8
$wsOff.Range("G3").Value = 8

# DEF sheet - row 3 ("R") gets updated with Wild Card round stats added in
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 448
$wsDef.Range("C3").Value = 319
$wsDef.Range("D3").Value = 102
$wsDef.Range("E3").Value = 54

$wb.Save()
